# Edit: "Added info on HTTP error codes to presentation"
# 1) Merge the two runs on slide 16 ("Many facilities and " + "helpers built in")
#    into a single run with text "Many facilities and helpers built in".
# 2) Append three new "Title and Content" slides (17, 18, 19) describing
#    HTTP status codes.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Fix slide 16 bullet text (merge two runs into one)
# ---------------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
$content16 = $slide16.Shapes.Item(2).TextFrame.TextRange
$bullet3 = $content16.Paragraphs(3)
# Force a real text change so the two runs actually get collapsed into one.
$bullet3.Text = "placeholder"
$bullet3.Text = "Many facilities and helpers built in"

# ---------------------------------------------------------------------------
# 2) New slide 17 - "HTTP Status Codes"
# ---------------------------------------------------------------------------
$s17 = $p.Slides.Add(17, 2)
$s17.Shapes.Item(1).TextFrame.TextRange.Text = "HTTP Status Codes"

$c17 = $s17.Shapes.Item(2).TextFrame.TextRange
$c17.Text = "1xx " + [char]0x2013 + " Informational " + [char]0x2013 + " seldom used`r" + `
            "2xx " + [char]0x2013 + " Success`r" + `
            "3xx " + [char]0x2013 + " Redirection, unchanged. Client should do something different to complete the request.`r" + `
            "4xx " + [char]0x2013 + " Client action caused an error`r" + `
            "5xx " + [char]0x2013 + " Server error"

# ---------------------------------------------------------------------------
# 3) New slide 18 - "Common HTTP Status Codes" (part 1)
# ---------------------------------------------------------------------------
$s18 = $p.Slides.Add(18, 2)
$s18.Shapes.Item(1).TextFrame.TextRange.Text = "Common HTTP Status Codes"

$c18 = $s18.Shapes.Item(2).TextFrame.TextRange
$c18.Text = "200 OK " + [char]0x2013 + " Request worked. Nothing to report`r" + `
            "201 Created " + [char]0x2013 + " Indicate new resource created`r" + `
            "204 No Content " + [char]0x2013 + " Used with conditional GETs`r" + `
            "`r" + `
            "304 Not Modified " + [char]0x2013 + " Used with `r"

$c18last = $c18.Paragraphs(6)
$c18last.ParagraphFormat.Bullet.Visible = 0
$c18last.Font.Size = 30

# ---------------------------------------------------------------------------
# 4) New slide 19 - "Common HTTP Status Codes" (part 2)
# ---------------------------------------------------------------------------
$s19 = $p.Slides.Add(19, 2)
$s19.Shapes.Item(1).TextFrame.TextRange.Text = "Common HTTP Status Codes"

$c19 = $s19.Shapes.Item(2).TextFrame.TextRange
$c19.Text = "401 Unauthorized`r" + `
            "404 Not found`r" + `
            "`r" + `
            "500 Internal Server Error " + [char]0x2013 + " Include details in body`r" + `
            "503 Service Unavailable " + [char]0x2013 + " Under maint, etc.`r"

# Split the last non-blank bullet into three runs: "...Under ", "maint", ", etc."
$bullet5 = $c19.Paragraphs(5)
$bullet5.Text = "503 Service Unavailable " + [char]0x2013 + " Under "
$null = $bullet5.InsertAfter("maint")
$null = $bullet5.InsertAfter(", etc.")
